# Update "想去人数" (interested-people count) values in column F across
# the four worksheets, per the upstream data refresh (gh-pages output
# generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        6  = 861
        7  = 449
        9  = 2110
        10 = 605
        11 = 269
        12 = 109
        13 = 1015
        14 = 161
        15 = 2159
        16 = 613
        17 = 11268
        18 = 1141
        19 = 546
    }
    "演出" = @{
        2  = 7
        9  = 145
        11 = 74
    }
    "本地生活" = @{
        4 = 444
    }
    "全部类型" = @{
        5  = 444
        6  = 7
        10 = 861
        12 = 449
        15 = 2110
        16 = 605
        17 = 269
        19 = 109
        21 = 1015
        23 = 161
        24 = 145
        26 = 2159
        27 = 613
        28 = 11268
        29 = 74
        31 = 1141
        32 = 546
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
